{"js": "const replacements = [\n  [\"N = 87,276\", \"N = 33,143\"],\n  [\"96.5 (89.8, 102.7)\", \"96.5 (89.8, 102.8)\"],\n  [\"462.2 (318.4, 638.1)\", \"465.2 (321.9, 637.8)\"],\n  [\"232.5 (114.0, 401.1)\", \"235.1 (116.5, 407.3)\"],\n  [\"584.2 (371.7, 849.1)\", \"587.5 (375.1, 853.5)\"],\n  [\"1,917 (2.2)\", \"738 (2.2)\"],\n  [\"1,227 (1.4)\", \"470 (1.4)\"],\n  [\"63.2 (56.1, 68.4)\", \"63.6 (56.5, 68.6)\"],\n  [\"2,590 (3.0)\", \"994 (3.0)\"],\n  [\"84,686 (97)\", \"32,149 (97)\"],\n  [\"50,371 (58)\", \"19,169 (58)\"],\n  [\"36,905 (42)\", \"13,974 (42)\"],\n  [\"6,683 (7.7)\", \"2,564 (7.7)\"],\n  [\"21,174 (24)\", \"7,956 (24)\"],\n  [\"20,530 (24)\", \"7,741 (23)\"],\n  [\"38,889 (45)\", \"14,882 (45)\"],\n  [\"11,034 (13)\", \"4,244 (13)\"],\n  [\"18,743 (21)\", \"7,236 (22)\"],\n  [\"22,949 (26)\", \"8,680 (26)\"],\n  [\"20,512 (24)\", \"7,674 (23)\"],\n  [\"6,000 (6.9)\", \"2,221 (6.7)\"],\n  [\"8,038 (9.2)\", \"3,088 (9.3)\"],\n  [\"-2.5 (-3.8, -0.2)\", \"-2.5 (-3.9, -0.3)\"],\n  [\"50,562 (58)\", \"19,136 (58)\"],\n  [\"30,886 (35)\", \"11,801 (36)\"],\n  [\"5,828 (6.7)\", \"2,206 (6.7)\"],\n  [\"4,765 (5.5)\", \"1,819 (5.5)\"],\n  [\"17,618 (20)\", \"6,596 (20)\"],\n  [\"21,962 (25)\", \"8,187 (25)\"],\n  [\"22,953 (26)\", \"8,744 (26)\"],\n  [\"19,978 (23)\", \"7,797 (24)\"],\n  [\"62,223 (71)\", \"23,637 (71)\"],\n  [\"21,985 (25)\", \"8,369 (25)\"],\n  [\"3,068 (3.5)\", \"1,137 (3.4)\"],\n  [\"15,701 (18)\", \"5,976 (18)\"],\n  [\"29,625 (34)\", \"11,299 (34)\"],\n  [\"41,950 (48)\", \"15,868 (48)\"],\n  [\"72,160 (83)\", \"27,412 (83)\"],\n  [\"14,426 (17)\", \"5,477 (17)\"],\n  [\"690 (0.8)\", \"254 (0.8)\"],\n  [\"74,100 (85)\", \"28,209 (85)\"],\n  [\"12,770 (15)\", \"4,786 (14)\"],\n  [\"406 (0.5)\", \"148 (0.4)\"],\n  [\"28,516 (33)\", \"10,769 (32)\"],\n  [\"52,300 (60)\", \"19,946 (60)\"],\n  [\"6,460 (7.4)\", \"2,428 (7.3)\"],\n  [\"14,915 (17)\", \"5,594 (17)\"],\n  [\"64,378 (74)\", \"24,492 (74)\"],\n  [\"7,983 (9.1)\", \"3,057 (9.2)\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Renamed-variable numeric updates for the PA1/PA3 table (testTable2_pa4).\n# Each table cell holds exactly one paragraph whose whole visible text is one\n# of the old statistic strings below; replace it in place with the new one,\n# preserving the run's formatting (we only touch the text node of each match).\n$d = $word.ActiveDocument\n\n$replacements = @{}\n$replacements[\"N = 87,276\"] = \"N = 33,143\"\n$replacements[\"96.5 (89.8, 102.7)\"] = \"96.5 (89.8, 102.8)\"\n$replacements[\"462.2 (318.4, 638.1)\"] = \"465.2 (321.9, 637.8)\"\n$replacements[\"232.5 (114.0, 401.1)\"] = \"235.1 (116.5, 407.3)\"\n$replacements[\"584.2 (371.7, 849.1)\"] = \"587.5 (375.1, 853.5)\"\n$replacements[\"1,917 (2.2)\"] = \"738 (2.2)\"\n$replacements[\"1,227 (1.4)\"] = \"470 (1.4)\"\n$replacements[\"63.2 (56.1, 68.4)\"] = \"63.6 (56.5, 68.6)\"\n$replacements[\"2,590 (3.0)\"] = \"994 (3.0)\"\n$replacements[\"84,686 (97)\"] = \"32,149 (97)\"\n$replacements[\"50,371 (58)\"] = \"19,169 (58)\"\n$replacements[\"36,905 (42)\"] = \"13,974 (42)\"\n$replacements[\"6,683 (7.7)\"] = \"2,564 (7.7)\"\n$replacements[\"21,174 (24)\"] = \"7,956 (24)\"\n$replacements[\"20,530 (24)\"] = \"7,741 (23)\"\n$replacements[\"38,889 (45)\"] = \"14,882 (45)\"\n$replacements[\"11,034 (13)\"] = \"4,244 (13)\"\n$replacements[\"18,743 (21)\"] = \"7,236 (22)\"\n$replacements[\"22,949 (26)\"] = \"8,680 (26)\"\n$replacements[\"20,512 (24)\"] = \"7,674 (23)\"\n$replacements[\"6,000 (6.9)\"] = \"2,221 (6.7)\"\n$replacements[\"8,038 (9.2)\"] = \"3,088 (9.3)\"\n$replacements[\"-2.5 (-3.8, -0.2)\"] = \"-2.5 (-3.9, -0.3)\"\n$replacements[\"50,562 (58)\"] = \"19,136 (58)\"\n$replacements[\"30,886 (35)\"] = \"11,801 (36)\"\n$replacements[\"5,828 (6.7)\"] = \"2,206 (6.7)\"\n$replacements[\"4,765 (5.5)\"] = \"1,819 (5.5)\"\n$replacements[\"17,618 (20)\"] = \"6,596 (20)\"\n$replacements[\"21,962 (25)\"] = \"8,187 (25)\"\n$replacements[\"22,953 (26)\"] = \"8,744 (26)\"\n$replacements[\"19,978 (23)\"] = \"7,797 (24)\"\n$replacements[\"62,223 (71)\"] = \"23,637 (71)\"\n$replacements[\"21,985 (25)\"] = \"8,369 (25)\"\n$replacements[\"3,068 (3.5)\"] = \"1,137 (3.4)\"\n$replacements[\"15,701 (18)\"] = \"5,976 (18)\"\n$replacements[\"29,625 (34)\"] = \"11,299 (34)\"\n$replacements[\"41,950 (48)\"] = \"15,868 (48)\"\n$replacements[\"72,160 (83)\"] = \"27,412 (83)\"\n$replacements[\"14,426 (17)\"] = \"5,477 (17)\"\n$replacements[\"690 (0.8)\"] = \"254 (0.8)\"\n$replacements[\"74,100 (85)\"] = \"28,209 (85)\"\n$replacements[\"12,770 (15)\"] = \"4,786 (14)\"\n$replacements[\"406 (0.5)\"] = \"148 (0.4)\"\n$replacements[\"28,516 (33)\"] = \"10,769 (32)\"\n$replacements[\"52,300 (60)\"] = \"19,946 (60)\"\n$replacements[\"6,460 (7.4)\"] = \"2,428 (7.3)\"\n$replacements[\"14,915 (17)\"] = \"5,594 (17)\"\n$replacements[\"64,378 (74)\"] = \"24,492 (74)\"\n$replacements[\"7,983 (9.1)\"] = \"3,057 (9.2)\"\n\n# Table-cell paragraphs end with a paragraph mark (CR) followed by the cell\n# mark (BEL); trim both so we compare/replace only the visible text.\n$cr = [char]13\n$cellMark = [char]7\n$appliedCount = 0\n\nforeach ($p in $d.Paragraphs) {\n    $visibleText = $p.Range.Text.TrimEnd($cr, $cellMark)\n    if ($replacements.ContainsKey($visibleText)) {\n        $p.Range.Text = $replacements[$visibleText]\n        $appliedCount++\n    }\n}\n\nif ($appliedCount -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) replacements, only applied $appliedCount\"\n}\n"}
